$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.024.68'
$ws.Range('E2').Value = '  -1.98%  '
$ws.Range('D3').Value = '3.072.57'
$ws.Range('E3').Value = '  -1.38%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '521.13'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.29%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.47'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.79%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').Value = '3.073.27'
$ws.Range('E8').Value = '  -1.29%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.455'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.33'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.04%  '
$ws.Range('E11').Value = '  -2.71%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.396'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.86%  '
$ws.Range('D13').Value = '3.625.08'
$ws.Range('E13').Value = '  -0.72%  '
$ws.Range('E14').Value = '  +1.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.24'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.13%  '
$ws.Range('E16').Value = '  -2.89%  '
$ws.Range('D17').Value = '57.172.60'
$ws.Range('E17').Value = '  -1.78%  '
$ws.Range('D18').Value = '3.083.89'
$ws.Range('E18').Value = '  -1.26%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.86'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.43'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.82'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '347.53'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.21%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '68.58'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.39%  '
$ws.Range('E25').Value = '  -3.64%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.166'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.23%  '
$ws.Range('D28').Value = '0.0₃0864'
$ws.Range('E28').Value = '  -7.27%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.19'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.86'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.99%  '
$ws.Range('E32').Value = '  -8.35%  '
$ws.Range('E33').Value = '  -1.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.84'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.17%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '158.96'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.33%  '
$ws.Range('E36').Value = '  -5.15%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '25.53'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.23'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.60%  '
$ws.Range('E40').Value = '  -1.86%  '
$ws.Range('E41').Value = '  -4.16%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.02'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.43%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.691'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.19%  '
$ws.Range('D44').Value = '2.393.51'
$ws.Range('E44').Value = '  +5.32%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '36.63'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.48%  '
$ws.Range('E46').Value = '  +0.33%  '
$ws.Range('D47').Value = '3.131.55'
$ws.Range('E47').Value = '  -0.81%  '
$ws.Range('E48').Value = '  -0.40%  '
$ws.Range('E49').Value = '  -4.58%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.95'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.09%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.66'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.77%  '
